$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '96.952.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.676.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.71%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("E5").Value = '  -0.74%  '

$ws.Range("E6").Value = '  +9.84%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '653.75'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.37%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.423'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.59%  '

$ws.Range("E9").Value = '  +3.35%  '

$ws.Range("E10").Value = '  +0.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.673.04'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.44'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.59%  '

$ws.Range("E13").Value = '  +1.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.90'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.362.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.78%  '

$ws.Range("E16").Value = '  +3.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.740.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.03'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.697.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.530'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '530.84'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.13%  '

$ws.Range("E24").Value = '  +0.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.19%  '

$ws.Range("E26").Value = '  -0.90%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '102.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.62%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.872.11'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.66%  '

$ws.Range("E30").Value = '  -0.30%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.16%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.04'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.56%  '

$ws.Range("E33").Value = '  -0.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.89'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +15.31%  '

$ws.Range("E35").Value = '  +0.56%  '

$ws.Range("E36").Value = '  +0.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '32.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.37%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '655.67'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.52%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.605'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.13%  '

$ws.Range("E40").Value = '  +1.69%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.90'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +14.37%  '

$ws.Range("E42").Value = '  +5.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.68%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.962'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '38.19'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +15.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.452'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0460'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.55%  '

$ws.Range("E49").Value = '  +0.58%  '

$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.64'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.11%  '

$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.73'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.41%  '
